{"js": "const replacements = [\n  [\"2024-07-26 Friday\", \"2024-07-27 Saturday\"],\n  [\"525\u00d72=1050\", \"565\u00d76=3390\"],\n  [\"435\u00d77=3045\", \"547\u00d78=4376\"],\n  [\"373\u00d78=2984\", \"521\u00d74=2084\"],\n  [\"745\u00d78=5960\", \"819\u00d77=5733\"],\n  [\"490\u00d73=1470\", \"186\u00d72=372\"],\n  [\"503\u00d77=3521\", \"322\u00d77=2254\"],\n  [\"875\u00d79=7875\", \"797\u00d79=7173\"],\n  [\"593\u00d74=2372\", \"211\u00d79=1899\"],\n  [\"582\u00d79=5238\", \"889\u00d76=5334\"],\n  [\"227\u00d74=908\", \"154\u00d79=1386\"],\n  [\"946\u00d76=5676\", \"363\u00d74=1452\"],\n  [\"101\u00d77=707\", \"324\u00d77=2268\"],\n  [\"425\u00d77=2975\", \"933\u00d74=3732\"],\n  [\"451\u00d79=4059\", \"254\u00d73=762\"],\n  [\"482\u00d79=4338\", \"798\u00d76=4788\"],\n  [\"128\u00d75=640\", \"872\u00d76=5232\"],\n  [\"776\u00d75=3880\", \"521\u00d76=3126\"],\n  [\"889\u00d74=3556\", \"347\u00d77=2429\"],\n  [\"917\u00d77=6419\", \"671\u00d78=5368\"],\n  [\"665\u00d79=5985\", \"390\u00d76=2340\"],\n  [\"807\u00d75=4035\", \"988\u00d78=7904\"],\n  [\"987\u00d75=4935\", \"547\u00d78=4376\"],\n  [\"736\u00d76=4416\", \"987\u00d73=2961\"],\n  [\"247\u00d74=988\", \"179\u00d77=1253\"],\n  [\"689\u00d75=3445\", \"879\u00d74=3516\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-07-26 Friday\", \"2024-07-27 Saturday\"),\n    @(\"525\u00d72=1050\", \"565\u00d76=3390\"),\n    @(\"435\u00d77=3045\", \"547\u00d78=4376\"),\n    @(\"373\u00d78=2984\", \"521\u00d74=2084\"),\n    @(\"745\u00d78=5960\", \"819\u00d77=5733\"),\n    @(\"490\u00d73=1470\", \"186\u00d72=372\"),\n    @(\"503\u00d77=3521\", \"322\u00d77=2254\"),\n    @(\"875\u00d79=7875\", \"797\u00d79=7173\"),\n    @(\"593\u00d74=2372\", \"211\u00d79=1899\"),\n    @(\"582\u00d79=5238\", \"889\u00d76=5334\"),\n    @(\"227\u00d74=908\", \"154\u00d79=1386\"),\n    @(\"946\u00d76=5676\", \"363\u00d74=1452\"),\n    @(\"101\u00d77=707\", \"324\u00d77=2268\"),\n    @(\"425\u00d77=2975\", \"933\u00d74=3732\"),\n    @(\"451\u00d79=4059\", \"254\u00d73=762\"),\n    @(\"482\u00d79=4338\", \"798\u00d76=4788\"),\n    @(\"128\u00d75=640\", \"872\u00d76=5232\"),\n    @(\"776\u00d75=3880\", \"521\u00d76=3126\"),\n    @(\"889\u00d74=3556\", \"347\u00d77=2429\"),\n    @(\"917\u00d77=6419\", \"671\u00d78=5368\"),\n    @(\"665\u00d79=5985\", \"390\u00d76=2340\"),\n    @(\"807\u00d75=4035\", \"988\u00d78=7904\"),\n    @(\"987\u00d75=4935\", \"547\u00d78=4376\"),\n    @(\"736\u00d76=4416\", \"987\u00d73=2961\"),\n    @(\"247\u00d74=988\", \"179\u00d77=1253\"),\n    @(\"689\u00d75=3445\", \"879\u00d74=3516\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    [void]$range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
